$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row,A,B,C) appended to the sheet - C is blank for many rows
$newRowsText = @"
1921,36,71,1
1922,36,72,2
1923,36,73,3
1924,36,74,
1925,36,75,
1926,36,76,
1927,36,77,
1928,36,78,
1929,36,79,
1930,36,80,
1931,37,71,
1932,37,72,
1933,37,73,
1934,37,74,
1935,37,75,
1936,37,76,
1937,37,77,
1938,37,78,
1939,37,79,
1940,37,80,
1941,38,71,
1942,38,72,
1943,38,73,
1944,38,74,
1945,38,75,
1946,38,76,
1947,38,77,
1948,38,78,
1949,38,79,
1950,38,80,
1951,39,71,
1952,39,72,
1953,39,73,
1954,39,74,
1955,39,75,
1956,39,76,
1957,39,77,
1958,39,78,
1959,39,79,
1960,39,80,
1961,40,71,
1962,40,72,
1963,40,73,
1964,40,74,
1965,40,75,
1966,40,76,
1967,40,77,
1968,40,78,
1969,40,79,
1970,40,80,
1971,44,71,
1972,44,72,
1973,44,73,
1974,44,74,
1975,44,75,
1976,44,76,
1977,44,77,
1978,44,78,
1979,44,79,
1980,44,80,
1981,45,71,
1982,45,72,
1983,45,73,
1984,45,74,
1985,45,75,
1986,45,76,
1987,45,77,
1988,45,78,
1989,45,79,
1990,45,80,
1991,52,71,
1992,52,72,
1993,52,73,
1994,52,74,
1995,52,75,
1996,52,76,
1997,52,77,
1998,52,78,
1999,52,79,
2000,52,80,
"@

$lines = $newRowsText -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line.Split(",")
    $r = [int]$parts[0]
    $a = [int]$parts[1]
    $b = [int]$parts[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    if ($parts.Length -gt 3 -and $parts[3] -ne "") {
        $c = [int]$parts[3]
        $ws.Cells.Item($r, 3).Value = $c
    }
}

# Update view: active cell / selection to match the author's final cursor position
$ws.Activate() | Out-Null
$ws.Range("C1923").Select() | Out-Null

# Page setup (paper size / orientation) touched during this edit
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
